{"js": "// Update the benchmark numbers in the single-column results table.\n// Each row is a table cell containing one run of text (sometimes the\n// run holds several tab-separated values). We address cells positionally\n// (row index, column 0) and rewrite the cell's text range in place so the\n// existing run formatting (Times New Roman, sz 22) is preserved.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfunction setCellText(rowIndex, text) {\n  const cell = table.getCell(rowIndex, 0);\n  const range = cell.body.getRange();\n  range.insertText(text, Word.InsertLocation.replace);\n}\n\n// Top summary rows\nsetCellText(0, \"0M\");\nsetCellText(1, \"0M\");\nsetCellText(2, \"0M\");\nsetCellText(3, \"2424\");\n\n// Allocation-rate style rows further down\nsetCellText(5, \"0.03361\");\nsetCellText(6, \"0.00017\");\nsetCellText(7, \"0.00045\");\nsetCellText(8, \"0.00023\");\nsetCellText(9, \"0.00024\");\nsetCellText(10, \"0.00029\");\nsetCellText(11, \"0.41428\");\n\n// Trailing rows previously holding a full tab-separated line; collapse\n// each back down to the single leading value.\nsetCellText(43, \"99.86\");\nsetCellText(44, \"0.41\");\nsetCellText(45, \"287\");\n\nawait context.sync();\n", "ps1": "# Update the benchmark numbers in the single-column results table.\n# Each row is a table cell containing one run of text (sometimes the run\n# holds several tab-separated values). We address cells positionally via\n# Table.Cell(row, col) (1-based) and overwrite the cell Range's Text so the\n# existing run formatting (Times New Roman, sz 22) is preserved.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfunction Set-CellText($rowIndex1Based, $text) {\n    $cell = $tbl.Cell($rowIndex1Based, 1)\n    $cell.Range.Text = $text\n}\n\n# Top summary rows\nSet-CellText 1 \"0M\"\nSet-CellText 2 \"0M\"\nSet-CellText 3 \"0M\"\nSet-CellText 4 \"2424\"\n\n# Allocation-rate style rows further down\nSet-CellText 6 \"0.03361\"\nSet-CellText 7 \"0.00017\"\nSet-CellText 8 \"0.00045\"\nSet-CellText 9 \"0.00023\"\nSet-CellText 10 \"0.00024\"\nSet-CellText 11 \"0.00029\"\nSet-CellText 12 \"0.41428\"\n\n# Trailing rows previously holding a full tab-separated line; collapse\n# each back down to the single leading value.\nSet-CellText 44 \"99.86\"\nSet-CellText 45 \"0.41\"\nSet-CellText 46 \"287\"\n"}
